{"js": "// The document contains a table whose cells each hold a single arithmetic\n// equation (e.g. \"67-8=59\"). This edit swaps each old equation for the new\n// one, in the same left-to-right, top-to-bottom order as the source table.\n// Every old/new string below is unique across the whole document (verified\n// against the source XML), so an exact, case-sensitive search-and-replace of\n// the full cell text is unambiguous and cannot clobber an unrelated cell.\nconst replacements = [\n  [\"67-8=59\", \"67+9=76\"],\n  [\"33-0=33\", \"33+26=59\"],\n  [\"93-40=53\", \"91-63=28\"],\n  [\"42+17=59\", \"24+46=70\"],\n  [\"35+38=73\", \"94-9=85\"],\n  [\"76-33=43\", \"67+28=95\"],\n  [\"51-24=27\", \"21+32=53\"],\n  [\"53-5=48\", \"60-27=33\"],\n  [\"98-7=91\", \"82-52=30\"],\n  [\"45-12=33\", \"48-43=5\"],\n  [\"72-67=5\", \"16+12=28\"],\n  [\"96-5=91\", \"72+13=85\"],\n  [\"64-9=55\", \"14+48=62\"],\n  [\"10+30=40\", \"20+40=60\"],\n  [\"76-49=27\", \"41-26=15\"],\n  [\"42+35=77\", \"98-28=70\"],\n  [\"6+57=63\", \"76-29=47\"],\n  [\"90-8=82\", \"39+10=49\"],\n  [\"35+44=79\", \"49-46=3\"],\n  [\"58+32=90\", \"47-44=3\"],\n  [\"2+30=32\", \"95-44=51\"],\n  [\"29+29=58\", \"57-55=2\"],\n  [\"71+14=85\", \"10+21=31\"],\n  [\"99-27=72\", \"81+16=97\"],\n  [\"2+2=4\", \"29+34=63\"],\n  [\"75-27=48\", \"81-40=41\"],\n  [\"91-84=7\", \"21+74=95\"],\n  [\"43-7=36\", \"16-2=14\"],\n  [\"9+2=11\", \"43-14=29\"],\n  [\"85-69=16\", \"9-5=4\"],\n  [\"79-39=40\", \"6+28=34\"],\n  [\"5+74=79\", \"43-12=31\"],\n  [\"80+13=93\", \"96-3=93\"],\n  [\"72+23=95\", \"9+5=14\"],\n  [\"45+42=87\", \"79-54=25\"],\n  [\"81-36=45\", \"43+47=90\"],\n  [\"98-45=53\", \"2+71=73\"],\n  [\"15+60=75\", \"92-11=81\"],\n  [\"51+26=77\", \"0+58=58\"],\n  [\"96-40=56\", \"17+12=29\"],\n  [\"0+34=34\", \"22+41=63\"],\n  [\"8+75=83\", \"21+50=71\"],\n  [\"29+48=77\", \"77-25=52\"],\n  [\"60+2=62\", \"48-34=14\"],\n  [\"73-70=3\", \"66-26=40\"],\n  [\"13+7=20\", \"68-24=44\"],\n  [\"91-1=90\", \"15-14=1\"],\n  [\"6+49=55\", \"45-30=15\"],\n  [\"21+11=32\", \"20+49=69\"],\n  [\"24+31=55\", \"84-18=66\"],\n  [\"44+23=67\", \"20+17=37\"],\n  [\"61-1=60\", \"33+59=92\"],\n  [\"45+50=95\", \"31-25=6\"],\n  [\"81-27=54\", \"59-9=50\"],\n  [\"50-46=4\", \"81-12=69\"],\n  [\"86-26=60\", \"8+60=68\"],\n  [\"17+36=53\", \"24+30=54\"],\n  [\"70+2=72\", \"17+73=90\"],\n  [\"79-38=41\", \"47+19=66\"],\n  [\"31+62=93\", \"37+5=42\"],\n  [\"23+2=25\", \"92-28=64\"],\n  [\"69-27=42\", \"87+11=98\"],\n  [\"80-20=60\", \"90-47=43\"],\n  [\"67-56=11\", \"64-50=14\"],\n  [\"76-51=25\", \"12+85=97\"],\n  [\"16+33=49\", \"16+77=93\"],\n  [\"91-51=40\", \"60+22=82\"],\n  [\"29+24=53\", \"24+58=82\"],\n  [\"62+6=68\", \"42-20=22\"],\n  [\"84-36=48\", \"46-38=8\"],\n  [\"2+46=48\", \"86-74=12\"],\n  [\"10+76=86\", \"36-4=32\"],\n  [\"63-2=61\", \"95-52=43\"],\n  [\"53-20=33\", \"24+67=91\"],\n  [\"3+60=63\", \"33+1=34\"],\n  [\"89-7=82\", \"55+43=98\"],\n  [\"32+38=70\", \"50+24=74\"],\n  [\"37+51=88\", \"72+17=89\"],\n  [\"95-45=50\", \"85-71=14\"],\n  [\"18-11=7\", \"11+68=79\"],\n  [\"11+32=43\", \"72+27=99\"],\n  [\"31+38=69\", \"69-4=65\"],\n  [\"36+9=45\", \"5+14=19\"],\n  [\"52-24=28\", \"84+14=98\"],\n  [\"89-70=19\", \"97-30=67\"],\n  [\"58-31=27\", \"80-11=69\"],\n  [\"30-11=19\", \"0+53=53\"],\n  [\"21-7=14\", \"90-61=29\"],\n  [\"66+8=74\", \"13+28=41\"],\n  [\"49-31=18\", \"97-37=60\"],\n  [\"91+2=93\", \"22+33=55\"],\n  [\"49-13=36\", \"31+55=86\"],\n  [\"8+16=24\", \"25+61=86\"],\n  [\"34+14=48\", \"51-8=43\"],\n  [\"69-43=26\", \"62-39=23\"],\n  [\"73+21=94\", \"13+55=68\"],\n  [\"77-49=28\", \"30+0=30\"],\n  [\"95-83=12\", \"63+32=95\"],\n  [\"18-12=6\", \"1+85=86\"],\n  [\"60-55=5\", \"39-28=11\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n}\n", "ps1": "# The document contains a table whose cells each hold a single arithmetic\n# equation (e.g. \"67-8=59\"). This edit swaps each old equation for the new\n# one, in the same left-to-right, top-to-bottom order as the source table.\n# Every old/new string below is unique across the whole document (verified\n# against the source XML), so an exact, case-sensitive whole-document\n# Find/Replace of the full cell text is unambiguous and cannot clobber an\n# unrelated cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  ,@(\"67-8=59\", \"67+9=76\")\n  ,@(\"33-0=33\", \"33+26=59\")\n  ,@(\"93-40=53\", \"91-63=28\")\n  ,@(\"42+17=59\", \"24+46=70\")\n  ,@(\"35+38=73\", \"94-9=85\")\n  ,@(\"76-33=43\", \"67+28=95\")\n  ,@(\"51-24=27\", \"21+32=53\")\n  ,@(\"53-5=48\", \"60-27=33\")\n  ,@(\"98-7=91\", \"82-52=30\")\n  ,@(\"45-12=33\", \"48-43=5\")\n  ,@(\"72-67=5\", \"16+12=28\")\n  ,@(\"96-5=91\", \"72+13=85\")\n  ,@(\"64-9=55\", \"14+48=62\")\n  ,@(\"10+30=40\", \"20+40=60\")\n  ,@(\"76-49=27\", \"41-26=15\")\n  ,@(\"42+35=77\", \"98-28=70\")\n  ,@(\"6+57=63\", \"76-29=47\")\n  ,@(\"90-8=82\", \"39+10=49\")\n  ,@(\"35+44=79\", \"49-46=3\")\n  ,@(\"58+32=90\", \"47-44=3\")\n  ,@(\"2+30=32\", \"95-44=51\")\n  ,@(\"29+29=58\", \"57-55=2\")\n  ,@(\"71+14=85\", \"10+21=31\")\n  ,@(\"99-27=72\", \"81+16=97\")\n  ,@(\"2+2=4\", \"29+34=63\")\n  ,@(\"75-27=48\", \"81-40=41\")\n  ,@(\"91-84=7\", \"21+74=95\")\n  ,@(\"43-7=36\", \"16-2=14\")\n  ,@(\"9+2=11\", \"43-14=29\")\n  ,@(\"85-69=16\", \"9-5=4\")\n  ,@(\"79-39=40\", \"6+28=34\")\n  ,@(\"5+74=79\", \"43-12=31\")\n  ,@(\"80+13=93\", \"96-3=93\")\n  ,@(\"72+23=95\", \"9+5=14\")\n  ,@(\"45+42=87\", \"79-54=25\")\n  ,@(\"81-36=45\", \"43+47=90\")\n  ,@(\"98-45=53\", \"2+71=73\")\n  ,@(\"15+60=75\", \"92-11=81\")\n  ,@(\"51+26=77\", \"0+58=58\")\n  ,@(\"96-40=56\", \"17+12=29\")\n  ,@(\"0+34=34\", \"22+41=63\")\n  ,@(\"8+75=83\", \"21+50=71\")\n  ,@(\"29+48=77\", \"77-25=52\")\n  ,@(\"60+2=62\", \"48-34=14\")\n  ,@(\"73-70=3\", \"66-26=40\")\n  ,@(\"13+7=20\", \"68-24=44\")\n  ,@(\"91-1=90\", \"15-14=1\")\n  ,@(\"6+49=55\", \"45-30=15\")\n  ,@(\"21+11=32\", \"20+49=69\")\n  ,@(\"24+31=55\", \"84-18=66\")\n  ,@(\"44+23=67\", \"20+17=37\")\n  ,@(\"61-1=60\", \"33+59=92\")\n  ,@(\"45+50=95\", \"31-25=6\")\n  ,@(\"81-27=54\", \"59-9=50\")\n  ,@(\"50-46=4\", \"81-12=69\")\n  ,@(\"86-26=60\", \"8+60=68\")\n  ,@(\"17+36=53\", \"24+30=54\")\n  ,@(\"70+2=72\", \"17+73=90\")\n  ,@(\"79-38=41\", \"47+19=66\")\n  ,@(\"31+62=93\", \"37+5=42\")\n  ,@(\"23+2=25\", \"92-28=64\")\n  ,@(\"69-27=42\", \"87+11=98\")\n  ,@(\"80-20=60\", \"90-47=43\")\n  ,@(\"67-56=11\", \"64-50=14\")\n  ,@(\"76-51=25\", \"12+85=97\")\n  ,@(\"16+33=49\", \"16+77=93\")\n  ,@(\"91-51=40\", \"60+22=82\")\n  ,@(\"29+24=53\", \"24+58=82\")\n  ,@(\"62+6=68\", \"42-20=22\")\n  ,@(\"84-36=48\", \"46-38=8\")\n  ,@(\"2+46=48\", \"86-74=12\")\n  ,@(\"10+76=86\", \"36-4=32\")\n  ,@(\"63-2=61\", \"95-52=43\")\n  ,@(\"53-20=33\", \"24+67=91\")\n  ,@(\"3+60=63\", \"33+1=34\")\n  ,@(\"89-7=82\", \"55+43=98\")\n  ,@(\"32+38=70\", \"50+24=74\")\n  ,@(\"37+51=88\", \"72+17=89\")\n  ,@(\"95-45=50\", \"85-71=14\")\n  ,@(\"18-11=7\", \"11+68=79\")\n  ,@(\"11+32=43\", \"72+27=99\")\n  ,@(\"31+38=69\", \"69-4=65\")\n  ,@(\"36+9=45\", \"5+14=19\")\n  ,@(\"52-24=28\", \"84+14=98\")\n  ,@(\"89-70=19\", \"97-30=67\")\n  ,@(\"58-31=27\", \"80-11=69\")\n  ,@(\"30-11=19\", \"0+53=53\")\n  ,@(\"21-7=14\", \"90-61=29\")\n  ,@(\"66+8=74\", \"13+28=41\")\n  ,@(\"49-31=18\", \"97-37=60\")\n  ,@(\"91+2=93\", \"22+33=55\")\n  ,@(\"49-13=36\", \"31+55=86\")\n  ,@(\"8+16=24\", \"25+61=86\")\n  ,@(\"34+14=48\", \"51-8=43\")\n  ,@(\"69-43=26\", \"62-39=23\")\n  ,@(\"73+21=94\", \"13+55=68\")\n  ,@(\"77-49=28\", \"30+0=30\")\n  ,@(\"95-83=12\", \"63+32=95\")\n  ,@(\"18-12=6\", \"1+85=86\")\n  ,@(\"60-55=5\", \"39-28=11\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  # wdFindContinue=1, wdReplaceAll=2; MatchCase=$true, MatchWholeWord=$true,\n  # MatchWildcards=$false, Forward=$true, Wrap=wdFindContinue.\n  $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
